# Expanse_details.xlsx — "messy code and unclear logic" cleanup.
# Fix up existing rows (category casing, amounts, a re-dated rent entry)
# and append the rest of the expense log that was missing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: rent -> Rent, amount 12 -> 10000, date moved forward ---
$ws.Range("A2").Value = "Rent"
$ws.Range("B2").Value = 10000
$ws.Range("C2").Value = 45847.29180555556

# --- Row 3: garbage placeholder -> real "Family treat" entry ---
$ws.Range("A3").Value = "Family treat "
$ws.Range("B3").Value = 205000
$ws.Range("C3").Value = 45803.29180555556

# The date column (C) uses a custom date number format on the existing
# rows (C2/C3). Copy that formatting down onto the new date cells below
# so C4:C6 render/save with the same style instead of "General".
$ws.Range("C3").Copy() | Out-Null
$ws.Range("C4:C6").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Row 4: Snack for sister ---
$ws.Range("A4").Value = "Snack for sister "
$ws.Range("B4").Value = 500000
$ws.Range("C4").Value = 45798.29180555556

# --- Row 5: buy simon sinek's book ---
$ws.Range("A5").Value = "buy simon sinek's book"
$ws.Range("B5").Value = 100000
$ws.Range("C5").Value = 45798.29180555556

# --- Row 6: Service laptop ---
$ws.Range("A6").Value = "Service laptop "
$ws.Range("B6").Value = 100000
$ws.Range("C6").Value = 45768.29180555556
